$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(57).Insert()

$ws.Range("A57").Value = 4
$ws.Range("B57").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C57").Value = "Los Lagos"
$ws.Range("D57").Value = 44540
$ws.Range("E57").Value = 10
$ws.Range("F57").Value = 100112024
$ws.Range("G57").Value = "Choclo"
$ws.Range("H57").Value = "Dulce o Americano"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 300
$ws.Range("K57").Value = 24000
$ws.Range("L57").Value = 24000
$ws.Range("M57").Value = 24000
$ws.Range("N57").Value = "`$/malla 70 unidades"
$ws.Range("O57").Value = "Región de Arica y Parinacota"
$ws.Range("P57").Value = 343
$ws.Range("Q57").Value = 70
$ws.Range("R57").Value = "Hortaliza"
